# Swap the betting-odds data (columns B:AD) between specific pairs of rows.
# Column A (rank/index) stays untouched; only B..AD swap places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(71, 72),
    @(101, 102),
    @(109, 110),
    @(215, 216),
    @(263, 265),
    @(316, 317)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rangeA = $ws.Range("B$r1`:AD$r1")
    $rangeB = $ws.Range("B$r2`:AD$r2")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}
